$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record, matching the style of existing headers (A1:AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (W-L-T) for every player row
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD
    $ws.Cells.Item($r, 31).Value = 67   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
